$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Project " + "Vulpix" + " " -> single run "Project Vulpix "
#    (drops the spellStart/spellEnd proofErr pair around "Vulpix")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Project Vulpix ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Project Vulpix ", 1)

# ---------------------------------------------------------------------------
# 2) "This weekly summary report for the third week of Project Vulpix will"
#    -> "...for the fourth week of Project Vulpix will"
#    (drops the spellStart/spellEnd proofErr pair around "Vulpix")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("for the third week of Project Vulpix will", $true, $false, $false, $false, $false,
                         $true, 1, $false, "for the fourth week of Project Vulpix will", 1)

# ---------------------------------------------------------------------------
# 3) "having the players take turns" -> "having the player takes turns"
#    (drops the gramStart/gramEnd proofErr pair around "take")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("having the players take turns", $true, $false, $false, $false, $false,
                         $true, 1, $false, "having the player takes turns", 1)

# ---------------------------------------------------------------------------
# 4) "the actions a player" -> "the actions, a player"
#    (drops the gramStart/gramEnd proofErr pair around "actions")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("the actions a player", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the actions, a player", 1)

# ---------------------------------------------------------------------------
# 5) Remove the stray "_GoBack" bookmark that currently sits after
#    "Weekly Progress Report #4".
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 6) "cards to program but" -> "cards to program, but"
#    (drops the gramStart/gramEnd proofErr pair around "program")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("cards to program but", $true, $false, $false, $false, $false,
                         $true, 1, $false, "cards to program, but", 1)

# ---------------------------------------------------------------------------
# 7) Re-create the "_GoBack" bookmark right after "program," (before " but we").
# ---------------------------------------------------------------------------
$rngBm = $d.Content
$rngBm.Find.Execute("program, but", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$bmPos = $rngBm.Start + 8
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------------
# 8) Insert a brand-new paragraph about the GUI debugging tool right before
#    the "Next week we will be finishing up..." paragraph.
# ---------------------------------------------------------------------------
$guiAnchor = $d.Paragraphs.Item(8).Range
$guiAnchor.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(9)
$newPara.Range.Text = "We have also created a graphical user interface to serve as a debugging tool. We found that when debugging we were displaying too much information to the console when using print statements. We designed and coded a GUI that can display the information and show that our code is doing what it is supposed to do. We are still debating if this GUI should make its way into the final project. Currently, it was not intended to be anything but a debugging tool. If we finish early or get ahead it might be redesigned and incorporated into the final project. "

# ---------------------------------------------------------------------------
# 9) Append a new sentence about the AI portion of the project to the end of
#    the final paragraph.
# ---------------------------------------------------------------------------
$rngEnd = $d.Content
$rngEnd.Find.Execute("This makes attacks more tedious to program than other game functions.",
                      $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngEnd.Collapse(0)
$rngEnd.InsertAfter(" We also plan to start the AI portion of our project next week. We won’t have any code for the AI next week, but we should have the design of the AI started next week. ")
